# Add a "Save" column (H) to the s_vals sheet, mirroring the existing
# header style (copied from G1) and filling the data rows with 0 values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: set the text first, then copy the formatting from the
# neighboring header cell (G1) so it reuses the same cell style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
